$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("syn_data2")

# Insert a copy of syn_data2 right before "accidents" and rename it to T10I4D100K
$source.Copy($wb.Worksheets.Item("accidents"))
$wb.Worksheets.Item("syn_data2 (2)").Name = "T10I4D100K"

# Insert another copy of syn_data2 right before "accidents" and rename it to kosarak
$source.Copy($wb.Worksheets.Item("accidents"))
$wb.Worksheets.Item("syn_data2 (2)").Name = "kosarak"

# Make "kosarak" the active/selected sheet (activeTab index 4, 0-based)
$wb.Worksheets.Item("kosarak").Activate()
